$d = $word.ActiveDocument

# Helper: replace one paragraph's content by targeting its Range and
# calling InsertXML with a full <w:p>...</w:p> fragment (wrapped in the
# minimal OOXML package envelope Word expects). This lets us control the
# exact run/break/xml:space structure instead of relying on plain
# Find&Replace (which would lose the xml:space="preserve" hint on the
# surrounding text).
function Set-ParagraphXml([int]$index, [string]$paragraphInnerXml) {
    $para = $d.Paragraphs($index)
    $rng = $para.Range
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $paragraphInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg) | Out-Null
}

# --- Paragraph 1: heading title + URL -------------------------------------
$p1 = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r>' + `
    '<w:t>Review 129: [Short] Watch Your Steps: Local Image and Scene Editing by Text Instructions,  26.08.23</w:t>' + `
    '<w:br/>' + `
    '<w:t>https://arxiv.org/abs/2211.09800.pdf</w:t>' + `
    '</w:r></w:p>'
Set-ParagraphXml 1 $p1

# --- Paragraph 2: bold "Paper: ..." link -----------------------------------
$p2 = '<w:p><w:r><w:rPr><w:b/></w:rPr>' + `
    '<w:t>Paper: https://arxiv.org/abs/2308.08947v1</w:t>' + `
    '</w:r></w:p>'
Set-ParagraphXml 2 $p2

# --- Paragraph 5: big Hebrew review body -----------------------------------
$t1 = 'מודלי הדיפוזיה לגמרי השתלטו כמעט על כל המשימות של הראייה הממוחשבת. למשל עריכה של תמונות (למשל להחליף ציפור בפרפר)בהתאם לתיאור טקסטואלי כבר מזמן עושים רק באמצעות מודלי דיפוזיה חזקים כמו  InstructPix2Pix או IP2P בקצרה. למרות התוצאות המדהימות עדיין יש אי התאמות בין התמונה הערוכה לבין המקורית. '
$t2 = 'היום ב-#shorthebrewpapereviews נסקור מאמר שמנסה לתקן את אי דיוקים אלו בצורה די אלגנטית. בשלב הראשון המודל המוצע מאתר את מיקום הפיקסלים שאותם צריך לשנות(מסכה) ובשלב השני עורכים את התמונה רק באזורים של המסכה. כל זה נעשה באמצעות מודלי דיפוזיה באופן די אלגנטי. '
$t3 = 'בשלב הראשון מרעישים את התמונה המקורית (עד רמת רעש מסוימת שהיא מהווה הייפרפרמטר חשוב מאוד) משתמשים במודל IP2P כדי לשערך את הרעש נוסף עבור ללא תופסת טקסט לעריכה ויחד איתו. כלומר במקרה הראשון אנו מפעילים מודל דיפוזיה סטנדרטי (ללא עריכה) ובמקרה השני כן עורכים את התמונה בהתאם לתיאור הטקסטואלי. לאחר מכן מחשבים את הערך המוחלט של ההפרש בין השערוכים אלו, מקצצים את החריגים (עם IQR עם מקדם 1.5). המסכה מקבלת ערך 1 (פיקסלים לעריכה) במקומות שההפרש הזה עולה על סף מסוים (הייפרפרמטר נוסף). '
$t4 = 'בשלב השני מרעישים את התמונה (רמת הרעש עוד הייפרפרמטר). ואז באמצעות מסירים את הרעש עם מודל IP2P (עם תיאור טקסטואלי) באיזורים של המסכה ובכל האזורים האחרים עושים זאת עם מודל דיפוזיה רגיל (הטקסט המוסף הוא ריק). בנוסף המחברים מכלילים את הגישה שלהם ל-NeRF (ייצוג של מודלי 3D). בגדול עושים את מה שמתואר למעלה על views מכל הזווית תוך שמירה של קוהרנטיות ביניהם.'

$p5 = '<w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r>' + `
    '<w:t xml:space="preserve">' + $t1 + '</w:t><w:br/><w:br/>' + `
    '<w:t xml:space="preserve">' + $t2 + '</w:t><w:br/><w:br/>' + `
    '<w:t xml:space="preserve">' + $t3 + '</w:t><w:br/><w:br/>' + `
    '<w:t>' + $t4 + '</w:t>' + `
    '</w:r></w:p>'
Set-ParagraphXml 5 $p5
